$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 353
$ws1.Range("F4").Value = 1295
$ws1.Range("F10").Value = 3547
$ws1.Range("F13").Value = 73
$ws1.Range("F14").Value = 51
$ws1.Range("F16").Value = 615
$ws1.Range("F17").Value = 108
$ws1.Range("F18").Value = 773
$ws1.Range("F19").Value = 214
$ws1.Range("F21").Value = 59
$ws1.Range("F24").Value = 2755
$ws1.Range("F25").Value = 5242
$ws1.Range("F29").Value = 3094
$ws1.Range("F31").Value = 2273
$ws1.Range("F35").Value = 136
$ws1.Range("F37").Value = 316
$ws1.Range("F38").Value = 37
$ws1.Range("F39").Value = 469
$ws1.Range("F40").Value = 812
$ws1.Range("F42").Value = 7
$ws1.Range("F44").Value = 42

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 353
$ws4.Range("F4").Value = 1295
$ws4.Range("F10").Value = 3547
$ws4.Range("F13").Value = 73
$ws4.Range("F15").Value = 51
$ws4.Range("F17").Value = 615
$ws4.Range("F18").Value = 108
$ws4.Range("F19").Value = 773
$ws4.Range("F20").Value = 214
$ws4.Range("F22").Value = 59
$ws4.Range("F25").Value = 2755
$ws4.Range("F26").Value = 5242
$ws4.Range("F30").Value = 3094
$ws4.Range("F32").Value = 2273
$ws4.Range("F36").Value = 136
$ws4.Range("F38").Value = 316
$ws4.Range("F39").Value = 37
$ws4.Range("F40").Value = 469
$ws4.Range("F41").Value = 812
$ws4.Range("F43").Value = 7
$ws4.Range("F45").Value = 42
